# Swap the "Fecha", "Volumen", "Precio mínimo", "Precio máximo",
# "Precio promedio ponderado" and "Precio $/Kg" values between row 2 and
# row 4, and between row 3 and row 5 (columns D, M, N, O, P, S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $cols = @("D", "M", "N", "O", "P", "S")
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value()
        $valB = $rangeB.Value()
        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

Swap-Rows 2 4
Swap-Rows 3 5
